$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.637.39'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').Value = '2.580.80'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''518.51'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '''139.26'
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('D9').Value = '2.588.94'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = '''6.54'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').Value = '3.037.21'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').Value = '58.671.51'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = '''20.38'
$ws.Range('E16').Value = '  +1.29%  '
$ws.Range('D17').Value = '2.594.08'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '''0.0000132'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = '''337.93'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '''4.29'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '''10.13'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').Value = '''6.49'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '''66.07'
$ws.Range('E24').Value = '  +2.08%  '
$ws.Range('D25').Value = '''0.167'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = '''0.404'
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').Value = '''7.03'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '0.0₃0715'
$ws.Range('E30').Value = '  -4.04%  '
$ws.Range('D31').Value = '''5.94'
$ws.Range('E31').Value = '  -4.65%  '
$ws.Range('D32').Value = '''1.56'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').Value = '''18.72'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').Value = '''148.55'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '''3.95'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('D37').Value = '''36.32'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').Value = '''1.46'
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('D39').Value = '''0.823'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '''0.813'
$ws.Range('E40').Value = '  -3.21%  '
$ws.Range('D41').Value = '''3.49'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').Value = '''273.57'
$ws.Range('E43').Value = '  +2.06%  '
$ws.Range('D44').Value = '''10.76'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('D45').Value = '''0.588'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '''0.0949'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').Value = '''0.0519'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = '''18.46'
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('D49').Value = '1.983.65'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = '''0.0220'
$ws.Range('E50').Value = '  -0.07%  '
$ws.Range('D51').Value = '''4.48'
$ws.Range('E51').Value = '  -1.28%  '
